$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    ,@("LatencyTest", 1.0, 200.0, 0.0, 0.0, 43.0, 23.25581395348837, 8.795)
    ,@("LatencyTest", 1.0, 200.0, 0.0, 0.0, 32.0, 31.25, 6.515000000000001)
    ,@("LatencyTest", 1.0, 200.0, 0.0, 0.0, 26.0, 38.46153846153846, 5.2090000000000005)
    ,@("LatencyTest", 1.0, 200.0, 0.0, 0.0, 25.0, 40.0, 5.163)
    ,@("Loading Initialization", 10.0, 2000.0, 2000.0, 0.0, 31.0, 322.5806451612903, 6.596)
    ,@("Loading Initialization", 10.0, 2000.0, 2000.0, 0.0, 25.0, 400.0, 5.331)
    ,@("Loading Initialization", 10.0, 2000.0, 2000.0, 0.0, 26.0, 384.6153846153846, 5.609)
    ,@("Loading Server", 100.0, 20000.0, 20000.0, 0.0, 136.0, 735.2941176470588, 47.656)
    ,@("Loading Initialization", 10.0, 2000.0, 2000.0, 0.0, 28.0, 357.1428571428571, 5.862)
    ,@("Loading Server", 100.0, 20000.0, 20000.0, 0.0, 85.0, 1176.4705882352941, 35.322)
    ,@("Loading Initialization", 10.0, 2000.0, 2000.0, 0.0, 27.0, 370.3703703703704, 5.982)
    ,@("Loading Server", 100.0, 200000.0, 200000.0, 0.0, 176.0, 568.1818181818182, 373.163)
    ,@("Loading Initialization", 10.0, 2000.0, 2000.0, 0.0, 27.0, 370.3703703703704, 5.843)
    ,@("Loading Server", 100.0, 20000.0, 20000.0, 0.0, 83.0, 1204.8192771084337, 34.633)
    ,@("Loading Initialization", 10.0, 2000.0, 2000.0, 0.0, 27.0, 370.3703703703704, 5.613)
    ,@("Loading Initialization", 10.0, 2000.0, 2000.0, 0.0, 28.0, 357.1428571428571, 5.912)
    ,@("Loading Server", 100.0, 20000.0, 20000.0, 0.0, 116.0, 862.0689655172413, 43.914)
    ,@("Loading Initialization", 10.0, 2000.0, 2000.0, 0.0, 26.0, 384.6153846153846, 5.632)
    ,@("Loading Server", 100.0, 20000.0, 20000.0, 0.0, 144.0, 694.4444444444443, 51.174)
    ,@("Loading Initialization", 10.0, 2000.0, 2000.0, 0.0, 27.0, 370.3703703703704, 5.672)
    ,@("Loading Initialization", 10.0, 2000.0, 2000.0, 0.0, 26.0, 384.6153846153846, 5.589)
    ,@("Loading Server", 100.0, 20000.0, 20000.0, 0.0, 81.0, 1234.567901234568, 33.965)
)

$startRow = 19
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
